$d = $word.ActiveDocument

# 1) Drop the trailing period from the existing bullet sentence.
$r1 = $d.Content
$r1.Find.Execute(
    "Refactoring existing monolithic application into Spring Cloud Microservices.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Refactoring existing monolithic application into Spring Cloud Microservices", 2)

# 2) Type the new tail of the sentence right after it. Toggling Bold on/off
#    forces this text to stay its own run instead of re-merging with the
#    identically formatted run that precedes it.
$r2 = $d.Content
$r2.Find.Execute("Spring Cloud Microservices")
$r2.Collapse(0)
$r2.Text = " processing records asynchronously via Kafka."
$r2.Bold = 1
$r2.Bold = 0

# 3) Move the "_GoBack" bookmark (Word's "last edit" marker) so it sits right
#    after "...asynchronously" and before " via Kafka." — this also splits
#    that typed text into two runs around the bookmark, and adding a
#    bookmark named "_GoBack" automatically replaces the prior one.
$r3 = $d.Content
$r3.Find.Execute("asynchronously")
$r3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r3)

Write-Output "done"
